$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data - a basic Admin account
$ws.Range("A2").Value = 1
$ws.Range("C2").Value = "Admin"
$ws.Range("B2").Value = "Administrator"
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "IT"

# Column widths (best fit) to accommodate new content
$ws.Columns.Item(2).ColumnWidth = 11.666666666666666
$ws.Columns.Item(3).ColumnWidth = 5.666666666666667
$ws.Columns.Item(5).ColumnWidth = 5.498697916666667
$ws.Columns.Item(7).ColumnWidth = 10.166666666666666

# Move selection to E2
$ws.Range("E2").Select()
